$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "ANZ Standard"
$ws.Range("F1").Value = "ANZ Advanced"
$ws.Range("I1").Value = "CBA"
$ws.Range("L1").Value = "Westpac"
$ws.Range("O1").Value = "NAB"
$ws.Range("R1").Value = "St George"
$ws.Range("U1").Value = "Bankwest"
$ws.Range("X1").Value = "UBank Loyalty Bonus"
$ws.Range("AC1").Value = "Citibank"
$ws.Range("A2").Value = "2017-02-19 21:50:48.032549"
$ws.Range("B2").Value = "Short term"
$ws.Range("C2").Value = 90
$ws.Range("D2").Value = 2.1
$ws.Range("F2").Value = 90
$ws.Range("G2").Value = 2.1
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 2.05
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = 2.05
$ws.Range("O2").Value = 90
$ws.Range("P2").Value = 2.1
$ws.Range("R2").Value = 90
$ws.Range("S2").Value = 2.1
$ws.Range("U2").Value = 90
$ws.Range("V2").Value = 2.45
$ws.Range("X2").Value = 90
$ws.Range("Y2").Value = 2.1
$ws.Range("AC2").Value = 90
$ws.Range("B3").Value = "Mid term"
$ws.Range("C3").Value = 180
$ws.Range("D3").Value = 2.2
$ws.Range("F3").Value = 120
$ws.Range("G3").Value = 2.2
$ws.Range("I3").Value = 210
$ws.Range("J3").Value = 2.1
$ws.Range("L3").Value = 180
$ws.Range("M3").Value = 2.1
$ws.Range("O3").Value = 240
$ws.Range("P3").Value = 2.2
$ws.Range("R3").Value = 180
$ws.Range("S3").Value = 2.2
$ws.Range("U3").Value = 210
$ws.Range("V3").Value = 2.6
$ws.Range("X3").Value = 180
$ws.Range("Y3").Value = 2.25
$ws.Range("AC3").Value = 180
$ws.Range("B4").Value = "Long term"
$ws.Range("C4").Value = 360
$ws.Range("D4").Value = 2.55
$ws.Range("F4").Value = 360
$ws.Range("G4").Value = 2.55
$ws.Range("I4").Value = 360
$ws.Range("J4").Value = 2.35
$ws.Range("L4").Value = 360
$ws.Range("M4").Value = 2.35
$ws.Range("O4").Value = 360
$ws.Range("P4").Value = 2.4
$ws.Range("R4").Value = 360
$ws.Range("S4").Value = 2.55
$ws.Range("U4").Value = 360
$ws.Range("V4").Value = 2.65
$ws.Range("X4").Value = 360
$ws.Range("Y4").Value = 2.5
$ws.Range("AC4").Value = 360
